$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> [newD, newE] ($null means "leave unchanged")
$changes = @{
    2  = @("26.846.08", "  +0.39%  ")
    3  = @("1.640.36", "  -0.22%  ")
    4  = @($null, "  -0.10%  ")
    5  = @("217.78", "  +0.57%  ")
    6  = @("0.498", "  -0.76%  ")
    7  = @($null, "  -0.15%  ")
    8  = @($null, "  -0.51%  ")
    9  = @($null, "  -0.81%  ")
    10 = @("19.24", "  +0.70%  ")
    11 = @("0.0844", "  +0.19%  ")
    12 = @("1.869.42", "  -0.03%  ")
    13 = @("1.641.67", "  -0.02%  ")
    14 = @($null, "  -0.65%  ")
    15 = @("0.527", "  +0.05%  ")
    16 = @("65.22", "  +1.22%  ")
    17 = @("26.835.48", "  +0.34%  ")
    18 = @("0.0₃0730", "  -0.67%  ")
    19 = @("215.27", "  +0.22%  ")
    20 = @($null, "  -0.21%  ")
    21 = @($null, "  -0.22%  ")
    22 = @("6.56", "  +5.02%  ")
    23 = @("2.39", "  -1.50%  ")
    24 = @("9.20", "  -1.53%  ")
    25 = @("147.22", "  +1.15%  ")
    26 = @($null, "  -0.28%  ")
    27 = @($null, "  -0.51%  ")
    28 = @("7.20", "  +1.24%  ")
    29 = @("15.73", "  +0.38%  ")
    30 = @($null, "  -0.26%  ")
    31 = @($null, "  +1.03%  ")
    32 = @($null, "  +0.38%  ")
    34 = @($null, "  +1.28%  ")
    35 = @("1.272.28", "  -1.44%  ")
    36 = @($null, "  +0.13%  ")
    37 = @($null, "  -1.84%  ")
    38 = @("0.530", "  -0.95%  ")
    39 = @("0.819", "  +0.06%  ")
    40 = @($null, "  -0.18%  ")
    41 = @("0.803", "  -0.47%  ")
    42 = @("5.32", "  -0.19%  ")
    43 = @("1.779.56", "  -0.58%  ")
    44 = @("92.41", "  +1.07%  ")
    45 = @("60.94", "  -0.88%  ")
    46 = @("2.05", "  -8.35%  ")
    47 = @("1.60", "  -0.49%  ")
    48 = @($null, "  -2.15%  ")
    51 = @($null, "  +0.04%  ")
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $newD = $vals[0]
    $newE = $vals[1]
    if ($null -ne $newD) {
        # Column D values are numeric-looking text (e.g. "9.20", "26.846.08")
        # and must stay plain text, matching the source data - force the
        # cell's number format to Text before writing so Excel doesn't
        # silently reinterpret/round them as numbers.
        $ws.Cells.Item($row, 4).NumberFormat = "@"
        $ws.Cells.Item($row, 4).Value = $newD
    }
    if ($null -ne $newE) {
        $ws.Cells.Item($row, 5).Value = $newE
    }
}

# Rows 49/50: EnergySwap and Algorand swap positions, each carrying new D/E values.
$ws.Cells.Item(49, 2).Value = "Algorand"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0967"
$ws.Cells.Item(49, 5).Value = "  -0.64%  "

$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "7.56"
$ws.Cells.Item(50, 5).Value = "  -1.34%  "
